$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "JP_acc"
$ws.Range("A3").Value = "JP_loose"
$ws.Range("A6").Value = "KR_loose"
$ws.Range("G7").Value = ""

$ws.Range("A5").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B6").Select()
